$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.403.77'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '2.599.49'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.59'
$ws.Range("E5").Value = '  +2.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.12'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.48'
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.334'
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").Value = '3.057.13'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '59.326.37'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.67'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.583.84'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '340.56'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.36'
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.07'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.35'
$ws.Range("E21").Value = '  -2.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.36'
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.164'
$ws.Range("E25").Value = '  -1.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.21'
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("D28").Value = '0.0₃0742'
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +4.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.81'
$ws.Range("E31").Value = '  -1.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.77'
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.78'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.96'
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.840'
$ws.Range("E36").Value = '  +2.86%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.821'
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '270.76'
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.598'
$ws.Range("E42").Value = '  +1.45%  '
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0952'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0524'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.54'
$ws.Range("E46").Value = '  +3.11%  '
$ws.Range("D47").Value = '1.939.34'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0222'
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.21'
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("E51").Value = '  +0.97%  '
